# Applies the "update a lot of thing" edit to Single-Level.xlsx
# - Row 28: clear C28 (was "Block Pointer"); D28 keeps "The ordering key field"
# - Row 29: C29 header text changes from "Block Pointer (P)" placeholder wording -> stays "Block Pointer (P)" (no-op, value already correct)
# - Row 33: fill in the previously-empty row with real data (40, 4096, 6, 9, 600000)
# - Row 45: rename headers, add new column D header "A Linear search on the data file = b/2"
# - Rows 46-55: add new column D formulas (=B{row}/2)
# - Row 58: remove column F header, rename remaining headers to the "binary/linear search" wording
# - Rows 59-68: remove column F (old "Linear search" values), shared-formula bookkeeping is
#   handled automatically by the engine when the column is deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 28 / 29 header cleanup
# ---------------------------------------------------------------------------
$ws.Range("C28").ClearContents()

# ---------------------------------------------------------------------------
# Row 33 now holds real data instead of being blank
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = 40
$ws.Range("B33").Value = 4096
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 9
$ws.Range("E33").Value = 600000

# ---------------------------------------------------------------------------
# Row 45 headers: reword C45, add new D45 header (copy the row's header style)
# ---------------------------------------------------------------------------
$ws.Range("C45").Value = "A binary search on the data file = [log2b]+"

$ws.Range("C45").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D45").Value = "A Linear search on the data file = b/2"

# ---------------------------------------------------------------------------
# Rows 46-55: new column D = B{row}/2 (copy the row's data style first)
# ---------------------------------------------------------------------------
$ws.Range("A46").Copy()
$ws.Range("D46:D55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D46").Formula = "=B46/2"
$ws.Range("D47").Formula = "=B47/2"
$ws.Range("D48").Formula = "=B48/2"
$ws.Range("D49").Formula = "=B49/2"
$ws.Range("D50").Formula = "=B50/2"
$ws.Range("D51").Formula = "=B51/2"
$ws.Range("D52").Formula = "=B52/2"
$ws.Range("D53").Formula = "=B53/2"
$ws.Range("D54").Formula = "=B54/2"
$ws.Range("D55").Formula = "=B55/2"

# ---------------------------------------------------------------------------
# Row 58 headers: reword remaining A:E headers
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "The size of each index Entry (Ri = Vssn + P)"
$ws.Range("B58").Value = "The blocking factor for the index (bfri = [B/Ri]-)"
$ws.Range("C58").Value = "The total number of index entries (ri = b)"
$ws.Range("D58").Value = "The number of index blocks is hence (bi = [ri/bfri]+)"
$ws.Range("E58").Value = "A binary search on the index file = [log2bi]+ +1"

# ---------------------------------------------------------------------------
# Remove column F entirely for rows 58-68 (old "Linear search = b/2" column,
# superseded by the new D46:D55 column above). Deleting shifts the shared
# formula bookkeeping / dimension automatically.
# ---------------------------------------------------------------------------
$ws.Range("F58:F68").Delete()

# ---------------------------------------------------------------------------
# Cosmetic view state: best-effort match of the recorded selection/scroll.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F32").Select()

Write-Output "edit applied"
